$d = $word.ActiveDocument

# --- 1. Text fixes (normalise stray spaces around punctuation) ---
# "Micaelo , en este documento estaré" -> "Micaelo, en este documento, estaré"
$d.Content.Find.Execute("Micaelo , en este documento estaré", $true, $false, $false, $false, $false, $true, 1, $false, "Micaelo, en este documento, estaré", 2)

# "proyecto . Pondré" -> "proyecto. Pondré"
$d.Content.Find.Execute("proyecto . Pondré", $true, $false, $false, $false, $false, $true, 1, $false, "proyecto. Pondré", 2)

# "repetida , así como documento ." -> "repetida así como documento."
$d.Content.Find.Execute("repetida , así como documento .", $true, $false, $false, $false, $false, $true, 1, $false, "repetida así como documento.", 2)

# --- 2. Re-seat the auto "_GoBack" bookmark where the last edit now sits: ---
#      right after "en este documento," and before " estaré"
$anchor = $d.Content
$anchor.Find.Execute("en este documento,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBack = $d.Range($anchor.End, $anchor.End)
$d.Bookmarks.Add("_GoBack", $goBack)

# --- 3. Paragraph formatting: no space after paragraph ---
$d.Paragraphs(1).SpaceAfter = 0

# --- 4. Three blank paragraphs appended at the end of the body ---
$tail = $d.Content
$blankPara = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$tail.Collapse(0)
$tail.InsertXML($blankPara)
$tail.Collapse(0)
$tail.InsertXML($blankPara)
$tail.Collapse(0)
$tail.InsertXML($blankPara)
